$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.926.89"
$ws.Range("E2").Value = "  -0.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.217.25"
$ws.Range("E3").Value = "  -1.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.24"
$ws.Range("E5").Value = "  +3.71%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  -0.52%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.08"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -4.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.36"
$ws.Range("E10").Value = "  +0.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  -2.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.90"
$ws.Range("E12").Value = "  -1.15%  "

# Row 13
$ws.Range("E13").Value = "  +0.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.553.18"
$ws.Range("E14").Value = "  -0.40%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.40"
$ws.Range("E15").Value = "  -1.61%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.216.49"
$ws.Range("E16").Value = "  -1.26%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("E17").Value = "  -3.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.862.93"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("E19").Value = "  -2.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.16"
$ws.Range("E20").Value = "  -0.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.39"
$ws.Range("E22").Value = "  -0.77%  "

# Row 23
$ws.Range("E23").Value = "  -1.50%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.33"
$ws.Range("E24").Value = "  -7.92%  "

# Row 25
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.60"
$ws.Range("E26").Value = "  -2.82%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.36"
$ws.Range("E27").Value = "  -1.58%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.70"
$ws.Range("E28").Value = "  +0.57%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.69%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  -2.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.45"
$ws.Range("E31").Value = "  +0.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.24"
$ws.Range("E32").Value = "  -0.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0850"
$ws.Range("E33").Value = "  +6.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.20"
$ws.Range("E34").Value = "  -2.91%  "

# Row 35
$ws.Range("E35").Value = "  -1.36%  "

# Row 36
$ws.Range("E36").Value = "  -0.97%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0353"
$ws.Range("E37").Value = "  +7.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.32"
$ws.Range("E38").Value = "  -0.54%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.42"
$ws.Range("E39").Value = "  -3.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.10"
$ws.Range("E40").Value = "  -1.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  +17.75%  "

# Row 42
$ws.Range("E42").Value = "  -3.34%  "

# Row 43
$ws.Range("E43").Value = "  -5.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "59.93"
$ws.Range("E44").Value = "  +0.09%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.69"
$ws.Range("E45").Value = "  -4.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0980"
$ws.Range("E47").Value = "  -1.11%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.454"
$ws.Range("E48").Value = "  -0.13%  "

# Row 49
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  -1.43%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.447.85"
$ws.Range("E51").Value = "  -0.04%  "
